# Apply edit to "Week 2" sheet, row 15 (task #9):
#  B15 = 10/02/2019 (date)
#  C15 = 17:20 (start time)
#  D15 = 18:15 (stop time)
#  G15 = "Prep."
#  H15 = "Reading HTML5 documentation and a little practice"
# then move selection to H16

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Week 2")

$ws.Range("B15").Value = 43506
$ws.Range("C15").Value = 0.72222222222222221
$ws.Range("D15").Value = 0.76041666666666663
$ws.Range("G15").Value = "Prep."
$ws.Range("H15").Value = "Reading HTML5 documentation and a little practice"

# The long comment now wraps to two lines in column H (wrap-text style),
# so the row grows to fit - mirror Excel's auto row-height here.
$ws.Rows.Item(15).RowHeight = 28.8

$ws.Activate()
$ws.Range("H16").Select()
